$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 565.3333
$ws.Range("I4").Value = 162.6
$ws.Range("J4").Value = 853
$ws.Range("K4").Value = 162.6
$ws.Range("L4").Value = 853
$ws.Range("M4").Value = -48.59999999999999
$ws.Range("N4").Value = -1081

$ws.Range("H18").Value = 6821.143
$ws.Range("I18").Value = 6437.5
$ws.Range("K18").Value = 6437.5
$ws.Range("M18").Value = -6153.5

$ws.Range("H33").Value = 912.4375
$ws.Range("I33").Value = 1019.1539
$ws.Range("K33").Value = 1019.1539
$ws.Range("M33").Value = -790.1539

$ws.Range("H86").Value = 5921.85
$ws.Range("I86").Value = 4461.385
$ws.Range("K86").Value = 4461.385
$ws.Range("M86").Value = -3338.385

$ws.Range("H89").Value = 5921.85
$ws.Range("I89").Value = 4461.385
$ws.Range("K89").Value = 22306.925
$ws.Range("M89").Value = -16690.925

$ws.Range("H99").Value = 1736.875
$ws.Range("I99").Value = 315.5
$ws.Range("J99").Value = 3158.25
$ws.Range("K99").Value = 946.5
$ws.Range("L99").Value = 9474.75
$ws.Range("M99").Value = 551.5
$ws.Range("N99").Value = -12470.75

$ws.Range("H112").Value = 5740.2173
$ws.Range("J112").Value = 5955.6816
$ws.Range("L112").Value = 17867.0448
$ws.Range("N112").Value = -20083.0448

$ws.Range("H131").Value = 5373.2666
$ws.Range("J131").Value = 10020.2
$ws.Range("L131").Value = 30060.6
$ws.Range("N131").Value = -40140.60000000001

$ws.Range("H138").Value = 3380.7778
$ws.Range("I138").Value = 1149.0526
$ws.Range("J138").Value = 4592.2856
$ws.Range("K138").Value = 3447.1578
$ws.Range("L138").Value = 13776.8568
$ws.Range("M138").Value = 1692.8422
$ws.Range("N138").Value = -24056.8568

$ws.Range("H141").Value = 12176.966
$ws.Range("I141").Value = 6338.3335
$ws.Range("J141").Value = 27503.375
$ws.Range("K141").Value = 19015.0005
$ws.Range("L141").Value = 82510.125
$ws.Range("M141").Value = -13835.0005
$ws.Range("N141").Value = -92870.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5563.274
$ws.Range("I32").Value = 3810.924
$ws.Range("K32").Value = 3810.924
$ws.Range("M32").Value = -3523.924

$ws.Range("H61").Value = 5071.6562
$ws.Range("I61").Value = 5138.516
$ws.Range("K61").Value = 5138.516
$ws.Range("M61").Value = -4926.516

$ws.Range("H122").Value = 1896341.2
$ws.Range("I122").Value = 2268
$ws.Range("K122").Value = 6804
$ws.Range("M122").Value = -4354

$ws.Range("H132").Value = 6451.1924
$ws.Range("I132").Value = 7374.8887
$ws.Range("J132").Value = 4372.875
$ws.Range("K132").Value = 22124.6661
$ws.Range("L132").Value = 13118.625
$ws.Range("M132").Value = -19594.6661
$ws.Range("N132").Value = -18178.625

$ws.Range("H136").Value = 5071.6562
$ws.Range("I136").Value = 5138.516
$ws.Range("K136").Value = 15415.548
$ws.Range("M136").Value = -12865.548

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 37650.5
$ws.Range("J92").Value = 37650.5
$ws.Range("L92").Value = 37650.5
$ws.Range("N92").Value = -42642.5

$ws.Range("H134").Value = 16578.834
$ws.Range("I134").Value = 18852.166
$ws.Range("K134").Value = 56556.49800000001
$ws.Range("M134").Value = -54021.49800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 10015
$ws.Range("J21").Value = 10015
$ws.Range("L21").Value = 10015
$ws.Range("N21").Value = -10485

$ws.Range("H58").Value = 2259.5557
$ws.Range("I58").Value = 1763.6786
$ws.Range("K58").Value = 1763.6786
$ws.Range("M58").Value = -1560.6786

$ws.Range("H122").Value = 1394.24
$ws.Range("I122").Value = 1146.1305
$ws.Range("J122").Value = 4247.5
$ws.Range("K122").Value = 3438.3915
$ws.Range("L122").Value = 12742.5
$ws.Range("M122").Value = -988.3914999999997
$ws.Range("N122").Value = -17642.5

$ws.Range("H132").Value = 120244.125
$ws.Range("I132").Value = 79310.16
$ws.Range("J132").Value = 297624.66
$ws.Range("K132").Value = 237930.48
$ws.Range("L132").Value = 892873.98
$ws.Range("M132").Value = -235400.48
$ws.Range("N132").Value = -897933.98

$ws.Range("H134").Value = 43946.965
$ws.Range("I134").Value = 73813.42999999999
$ws.Range("J134").Value = 11783.077
$ws.Range("K134").Value = 221440.29
$ws.Range("L134").Value = 35349.231
$ws.Range("M134").Value = -218905.29
$ws.Range("N134").Value = -40419.231

$ws.Range("H136").Value = 2259.5557
$ws.Range("I136").Value = 1763.6786
$ws.Range("K136").Value = 5291.0358
$ws.Range("M136").Value = -2741.0358

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 657.8889
$ws.Range("J68").Value = 604.2
$ws.Range("L68").Value = 1812.6
$ws.Range("N68").Value = -3434.6

$ws.Range("H71").Value = 657.8889
$ws.Range("J71").Value = 604.2
$ws.Range("L71").Value = 5437.8
$ws.Range("N71").Value = -13549.8

$ws.Range("H131").Value = 9693548
$ws.Range("I131").Value = 5557480.5
$ws.Range("J131").Value = 11909298
$ws.Range("K131").Value = 16672441.5
$ws.Range("L131").Value = 35727894
$ws.Range("M131").Value = -16667401.5
$ws.Range("N131").Value = -35737974

$ws.Range("H137").Value = 4032.4285
$ws.Range("J137").Value = 5332.6665
$ws.Range("L137").Value = 15997.9995
$ws.Range("N137").Value = -26197.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 508.375
$ws.Range("J13").Value = 508.375
$ws.Range("L13").Value = 508.375
$ws.Range("N13").Value = -786.375

$ws.Range("H41").Value = 4000
$ws.Range("I41").Value = 4000
$ws.Range("K41").Value = 4000
$ws.Range("M41").Value = -3645

$ws.Range("H80").Value = 2724899
$ws.Range("I80").Value = 8163564.5
$ws.Range("J80").Value = 5566.3335
$ws.Range("K80").Value = 8163564.5
$ws.Range("L80").Value = 5566.3335
$ws.Range("M80").Value = -8162566.5
$ws.Range("N80").Value = -7562.3335

$ws.Range("H83").Value = 2724899
$ws.Range("I83").Value = 8163564.5
$ws.Range("J83").Value = 5566.3335
$ws.Range("K83").Value = 40817822.5
$ws.Range("L83").Value = 27831.6675
$ws.Range("M83").Value = -40812830.5
$ws.Range("N83").Value = -37815.6675

$ws.Range("H102").Value = 2194250
$ws.Range("I102").Value = 2316412
$ws.Range("J102").Value = 1775408.6
$ws.Range("K102").Value = 2316412
$ws.Range("L102").Value = 1775408.6
$ws.Range("M102").Value = -2314790
$ws.Range("N102").Value = -1778652.6

$ws.Range("H122").Value = 623864.9399999999
$ws.Range("I122").Value = 972324.9399999999
$ws.Range("J122").Value = 76285
$ws.Range("K122").Value = 2916974.82
$ws.Range("L122").Value = 228855
$ws.Range("M122").Value = -2914524.82
$ws.Range("N122").Value = -233755

$ws.Range("H132").Value = 11070.77
$ws.Range("I132").Value = 8880.223
$ws.Range("K132").Value = 26640.669
$ws.Range("M132").Value = -24110.669

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3419.946
$ws.Range("I7").Value = 1642.3846
$ws.Range("J7").Value = 7621.4546
$ws.Range("K7").Value = 1642.3846
$ws.Range("L7").Value = 7621.4546
$ws.Range("M7").Value = -1530.3846
$ws.Range("N7").Value = -7845.4546

$ws.Range("H40").Value = 4161.5293
$ws.Range("I40").Value = 1775.7
$ws.Range("J40").Value = 7569.857
$ws.Range("K40").Value = 1775.7
$ws.Range("L40").Value = 7569.857
$ws.Range("M40").Value = -1639.7
$ws.Range("N40").Value = -7841.857

$ws.Range("H93").Value = 37037988
$ws.Range("I93").Value = 55556336
$ws.Range("J93").Value = 1296.6666
$ws.Range("K93").Value = 55556336
$ws.Range("L93").Value = 1296.6666
$ws.Range("M93").Value = -55555088
$ws.Range("N93").Value = -3792.6666

$ws.Range("H122").Value = 5208.696
$ws.Range("I122").Value = 3422.0908
$ws.Range("K122").Value = 10266.2724
$ws.Range("M122").Value = -7816.2724

$ws.Range("H126").Value = 3419.946
$ws.Range("I126").Value = 1642.3846
$ws.Range("J126").Value = 7621.4546
$ws.Range("K126").Value = 4927.1538
$ws.Range("L126").Value = 22864.3638
$ws.Range("M126").Value = -2457.1538
$ws.Range("N126").Value = -27804.3638

$ws.Range("H132").Value = 8050.5166
$ws.Range("I132").Value = 8962.653
$ws.Range("J132").Value = 3987.3635
$ws.Range("K132").Value = 26887.959
$ws.Range("L132").Value = 11962.0905
$ws.Range("M132").Value = -24357.959
$ws.Range("N132").Value = -17022.0905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 11503
$ws.Range("J31").Value = 12254.5
$ws.Range("L31").Value = 12254.5
$ws.Range("N31").Value = -12950.5

$ws.Range("H75").Value = 11000
$ws.Range("I75").Value = 11000
$ws.Range("K75").Value = 11000
$ws.Range("M75").Value = -10064

$ws.Range("H78").Value = 11000
$ws.Range("I78").Value = 11000
$ws.Range("K78").Value = 33000
$ws.Range("M78").Value = -28320

$ws.Range("H92").Value = 54549.5
$ws.Range("J92").Value = 54549.5
$ws.Range("L92").Value = 54549.5
$ws.Range("N92").Value = -59541.5

$ws.Range("H96").Value = 7749
$ws.Range("I96").Value = 7498
$ws.Range("K96").Value = 7498
$ws.Range("M96").Value = -6125

$ws.Range("H113").Value = 807.35297
$ws.Range("I113").Value = 529.61536
$ws.Range("K113").Value = 1588.84608
$ws.Range("M113").Value = 581.15392

$ws.Range("H122").Value = 3106.9
$ws.Range("I122").Value = 2206.2666
$ws.Range("K122").Value = 6618.7998
$ws.Range("M122").Value = -4168.7998

$ws.Range("H132").Value = 20227100
$ws.Range("I132").Value = 25649674
$ws.Range("K132").Value = 76949022
$ws.Range("M132").Value = -76946492
